# Adds the new "BARRANQUILLA" branch meeting-schedule rows (59-72) to the
# HorarioReuniones table, matching the upstream "Add files via upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- 1. Grow the table (and its ref / autofilter / dimension) by 14 rows ---
for ($i = 0; $i -lt 14; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# --- 2. Prepare the two reusable cell formats we need -------------------
# "plain" format: centered, light fill, no special number format
# (same look as the existing data rows, e.g. A11)
$ws.Range("A11").Copy()
$ws.Range("A59:E72").PasteSpecial(-4122)
$ws.Range("G59:G72").PasteSpecial(-4122)
$ws.Range("H59:H72").PasteSpecial(-4122)

# "time" format: centered, light fill, h:mm number format
# (same fill as A11/F2, but with the h:mm time format used by this sheet)
$ws.Range("F2").Copy()
$ws.Range("F59:F72").PasteSpecial(-4122)
$ws.Range("F59:F72").NumberFormat = "h:mm"

# H column time cells (all but H71, which stays blank)
$ws.Range("H59:H70").NumberFormat = "h:mm"
$ws.Range("H72").NumberFormat = "h:mm"

# --- 3. Row data, written left-to-right / top-to-bottom so new shared ---
#        strings are registered in the same order the sheet was authored.
#        ("German Vaquero" cells in column C are filled in last, below -
#        that matches how this workbook's string table was actually built.)
# Columns: A=SUCURSAL B=PRACTICANTE C=GERENTE D=PROYECTO
#          E=DIA INTERMEDIA F=HORA INTERMEDIA G=DIA SEMANAL H=HORA SEMANAL

$ws.Range("A59").Value = "BARRANQUILLA"
$ws.Range("B59").Value = "NA"
$ws.Range("C59").Value = "José Luis Suarez"
$ws.Range("D59").Value = "Amaretto"
$ws.Range("E59").Value = "Jueves"
$ws.Range("F59").Value = 0.375
$ws.Range("G59").Value = "Jueves"
$ws.Range("H59").Value = 0.58333333333333337

$ws.Range("A60").Value = "BARRANQUILLA"
$ws.Range("B60").Value = "NA"
$ws.Range("C60").Value = "José Luis Suarez"
$ws.Range("D60").Value = "Aquanova"
$ws.Range("E60").Value = "Miercoles"
$ws.Range("F60").Value = 0.375
$ws.Range("G60").Value = "Jueves"
$ws.Range("H60").Value = 0.35416666666666669

$ws.Range("A61").Value = "BARRANQUILLA"
$ws.Range("B61").Value = "NA"
$ws.Range("D61").Value = "CP Los Flamencos"
$ws.Range("E61").Value = "Jueves"
$ws.Range("F61").Value = 0.58333333333333337
$ws.Range("G61").Value = "Viernes"
$ws.Range("H61").Value = 0.58333333333333337

$ws.Range("A62").Value = "BARRANQUILLA"
$ws.Range("B62").Value = "NA"
$ws.Range("D62").Value = "CP Tayrona Aptos"
$ws.Range("E62").Value = "Jueves"
$ws.Range("F62").Value = 0.375
$ws.Range("G62").Value = "Viernes"
$ws.Range("H62").Value = 0.375

$ws.Range("A63").Value = "BARRANQUILLA"
$ws.Range("B63").Value = "NA"
$ws.Range("D63").Value = "CP Cocuy"
$ws.Range("E63").Value = "Martes"
$ws.Range("F63").Value = 0.35416666666666669
$ws.Range("G63").Value = "Miercoles"
$ws.Range("H63").Value = 0.35416666666666669

$ws.Range("A64").Value = "BARRANQUILLA"
$ws.Range("B64").Value = "NA"
$ws.Range("C64").Value = "José Luis Suarez"
$ws.Range("D64").Value = "Dimaro"
$ws.Range("E64").Value = "Martes"
$ws.Range("F64").Value = 0.375
$ws.Range("G64").Value = "Miercoles"
$ws.Range("H64").Value = 0.58333333333333337

$ws.Range("A65").Value = "BARRANQUILLA"
$ws.Range("B65").Value = "NA"
$ws.Range("C65").Value = "José Luis Suarez"
$ws.Range("D65").Value = "Genova"
$ws.Range("E65").Value = "Jueves"
$ws.Range("F65").Value = 0.41666666666666669
$ws.Range("G65").Value = "Viernes"
$ws.Range("H65").Value = 0.41666666666666669

$ws.Range("A66").Value = "BARRANQUILLA"
$ws.Range("B66").Value = "NA"
$ws.Range("D66").Value = "PD Arrecife"
$ws.Range("E66").Value = "Jueves"
$ws.Range("F66").Value = 0.375
$ws.Range("G66").Value = "Viernes"
$ws.Range("H66").Value = 0.375

$ws.Range("A67").Value = "BARRANQUILLA"
$ws.Range("B67").Value = "NA"
$ws.Range("D67").Value = "PD Cristalina"
$ws.Range("E67").Value = "Jueves"
$ws.Range("F67").Value = 0.35416666666666669
$ws.Range("G67").Value = "Viernes"
$ws.Range("H67").Value = 0.35416666666666669

$ws.Range("A68").Value = "BARRANQUILLA"
$ws.Range("B68").Value = "NA"
$ws.Range("D68").Value = "PD Marisma"
$ws.Range("E68").Value = "Jueves"
$ws.Range("F68").Value = 0.58333333333333337
$ws.Range("G68").Value = "Sabado"
$ws.Range("H68").Value = 0.3125

$ws.Range("A69").Value = "BARRANQUILLA"
$ws.Range("B69").Value = "NA"
$ws.Range("D69").Value = "PD Natura"
$ws.Range("E69").Value = "Jueves"
$ws.Range("F69").Value = 0.58333333333333337
$ws.Range("G69").Value = "Viernes"
$ws.Range("H69").Value = 0.39583333333333331

$ws.Range("A70").Value = "BARRANQUILLA"
$ws.Range("B70").Value = "NA"
$ws.Range("C70").Value = "José Luis Suarez"
$ws.Range("D70").Value = "Riverbay"
$ws.Range("E70").Value = "Martes"
$ws.Range("F70").Value = 0.375
$ws.Range("G70").Value = "Miercoles"
$ws.Range("H70").Value = 0.08333333333333333

$ws.Range("A71").Value = "BARRANQUILLA"
$ws.Range("B71").Value = "NA"
$ws.Range("C71").Value = "José Luis Suarez"
$ws.Range("D71").Value = "Riverblue"
$ws.Range("E71").Value = "Jueves"
$ws.Range("F71").Value = 0.375
$ws.Range("G71").Value = "-"
# H71 intentionally left blank (no recorded "HORA SEMANAL" time), but keeps
# its number/fill formatting from the PasteSpecial pass above.

$ws.Range("A72").Value = "BARRANQUILLA"
$ws.Range("B72").Value = "NA"
$ws.Range("D72").Value = "Solario"
$ws.Range("E72").Value = "Lunes"
$ws.Range("F72").Value = 0.33333333333333331
$ws.Range("G72").Value = "Martes"
$ws.Range("H72").Value = 0.33333333333333331

# --- 4. GERENTE column (C) - "German Vaquero" rows, filled last --------
$ws.Range("C61").Value = "Germán Vaquero"
$ws.Range("C62").Value = "Germán Vaquero"
$ws.Range("C63").Value = "Germán Vaquero"
$ws.Range("C66").Value = "Germán Vaquero"
$ws.Range("C67").Value = "Germán Vaquero"
$ws.Range("C68").Value = "Germán Vaquero"
$ws.Range("C69").Value = "Germán Vaquero"
$ws.Range("C72").Value = "Germán Vaquero"
